# Updates crypto price/volume data to match the latest scrape.
# Row 12/13 additionally swap rank positions: TRON moves up to row 12,
# WrappedEther moves down to row 13 (each keeps its own refreshed price/volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''27.861.04'
$ws.Range('E2').Value = '  -1.01%  '

# Row 3
$ws.Range('D3').Value = '''1.903.94'
$ws.Range('E3').Value = '  -0.59%  '

# Row 4
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.59%  '

# Row 5
$ws.Range('D5').Value = '''313.10'
$ws.Range('E5').Value = '  -1.29%  '

# Row 6
$ws.Range('E6').Value = '  -0.52%  '

# Row 7
$ws.Range('D7').Value = '''0.4972'
$ws.Range('E7').Value = '  +2.60%  '

# Row 8
$ws.Range('D8').Value = '''0.3821'
$ws.Range('E8').Value = '  -0.15%  '

# Row 9
$ws.Range('D9').Value = '''0.07338'
$ws.Range('E9').Value = '  -0.51%  '

# Row 10
$ws.Range('D10').Value = '''0.9101'
$ws.Range('E10').Value = '  -2.91%  '

# Row 11
$ws.Range('E11').Value = '  +0.34%  '

# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.07636'
$ws.Range('E12').Value = '  -2.30%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '''1.892.56'
$ws.Range('E13').Value = '  -1.29%  '

# Row 14
$ws.Range('D14').Value = '''5.488'
$ws.Range('E14').Value = '  -0.25%  '

# Row 15
$ws.Range('D15').Value = '''6.638'
$ws.Range('E15').Value = '  -0.02%  '

# Row 16
$ws.Range('D16').Value = '''91.38'
$ws.Range('E16').Value = '  +0.17%  '

# Row 17
$ws.Range('E17').Value = '  -0.71%  '

# Row 18
$ws.Range('D18').Value = '''0.000008731'
$ws.Range('E18').Value = '  -1.20%  '

# Row 19
$ws.Range('E19').Value = '  -0.43%  '

# Row 20
$ws.Range('D20').Value = '''27.890.19'
$ws.Range('E20').Value = '  -1.04%  '

# Row 21
$ws.Range('D21').Value = '''14.53'
$ws.Range('E21').Value = '  -2.22%  '

# Row 22
$ws.Range('D22').Value = '''5.134'
$ws.Range('E22').Value = '  -0.55%  '

# Row 23
$ws.Range('D23').Value = '''10.80'
$ws.Range('E23').Value = '  -1.08%  '

# Row 24
$ws.Range('D24').Value = '''154.48'
$ws.Range('E24').Value = '  -1.35%  '

# Row 25
$ws.Range('D25').Value = '''1.867'
$ws.Range('E25').Value = '  -2.74%  '

# Row 26
$ws.Range('D26').Value = '''2.232'
$ws.Range('E26').Value = '  +6.12%  '

# Row 27
$ws.Range('E27').Value = '  -0.79%  '

# Row 28
$ws.Range('D28').Value = '''115.24'

# Row 29
$ws.Range('D29').Value = '''4.950'
$ws.Range('E29').Value = '  -0.33%  '

# Row 30
$ws.Range('D30').Value = '''0.08973'

# Row 31
$ws.Range('D31').Value = '''3.194'
$ws.Range('E31').Value = '  -5.04%  '

# Row 32
$ws.Range('D32').Value = '''1.240'
$ws.Range('E32').Value = '  -1.01%  '

# Row 33
$ws.Range('D33').Value = '''0.7726'
$ws.Range('E33').Value = '  +0.12%  '

# Row 34
$ws.Range('D34').Value = '''4.643'
$ws.Range('E34').Value = '  -1.10%  '

# Row 35
$ws.Range('E35').Value = '  +0.76%  '

# Row 36
$ws.Range('D36').Value = '''2.577'
$ws.Range('E36').Value = '  -2.83%  '

# Row 37
$ws.Range('E37').Value = '  -0.15%  '

# Row 38
$ws.Range('D38').Value = '''0.5527'
$ws.Range('E38').Value = '  -0.03%  '

# Row 39
$ws.Range('E39').Value = '  -0.53%  '

# Row 40
$ws.Range('D40').Value = '''3.010'
$ws.Range('E40').Value = '  +0.24%  '

# Row 41
$ws.Range('D41').Value = '''6.994'
$ws.Range('E41').Value = '  -0.90%  '

# Row 42
$ws.Range('D42').Value = '''8.559'
$ws.Range('E42').Value = '  +1.14%  '

# Row 43
$ws.Range('D43').Value = '''0.1524'
$ws.Range('E43').Value = '  -0.32%  '

# Row 44
$ws.Range('D44').Value = '''111.56'
$ws.Range('E44').Value = '  +3.93%  '

# Row 45
$ws.Range('D45').Value = '''10.65'
$ws.Range('E45').Value = '  -0.76%  '

# Row 46
$ws.Range('D46').Value = '''0.4798'
$ws.Range('E46').Value = '  -1.30%  '

# Row 47
$ws.Range('D47').Value = '''0.9998'

# Row 48
$ws.Range('D48').Value = '''1.640'
$ws.Range('E48').Value = '  -1.07%  '

# Row 49
$ws.Range('D49').Value = '''67.53'
$ws.Range('E49').Value = '  -1.73%  '

# Row 50
$ws.Range('D50').Value = '''0.06077'
$ws.Range('E50').Value = '  -0.51%  '

# Row 51
$ws.Range('D51').Value = '''0.9004'
$ws.Range('E51').Value = '  -0.90%  '
